# Auto-generated edit script
# Applies numeric cell updates to the Kujata_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 11111991
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 18519318
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 18519318
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -18519456
$ws.Range("H116").Value = 3010.6843
$ws.Range("I116").Value = 2499.182
$ws.Range("J116").Value = 3714
$ws.Range("K116").Value = 2499.182
$ws.Range("L116").Value = 3714
$ws.Range("M116").Value = 942.8180000000002
$ws.Range("N116").Value = -10598
$ws.Range("H132").Value = 9665.5625
$ws.Range("I132").Value = 6614.9
$ws.Range("J132").Value = 14750
$ws.Range("K132").Value = 19844.7
$ws.Range("L132").Value = 44250
$ws.Range("M132").Value = -17314.7
$ws.Range("N132").Value = -49310
$ws.Range("H137").Value = 1918.4615
$ws.Range("I137").Value = 1309.25
$ws.Range("J137").Value = 2189.2222
$ws.Range("K137").Value = 3927.75
$ws.Range("L137").Value = 6567.6666
$ws.Range("M137").Value = -1377.75
$ws.Range("N137").Value = -11667.6666
$ws.Range("H138").Value = 2062.8223
$ws.Range("I138").Value = 831.75
$ws.Range("J138").Value = 2252.218
$ws.Range("K138").Value = 2495.25
$ws.Range("L138").Value = 6756.654
$ws.Range("M138").Value = 2644.75
$ws.Range("N138").Value = -17036.654

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3247.3333
$ws.Range("I32").Value = 3118.6223
$ws.Range("K32").Value = 3118.6223
$ws.Range("M32").Value = -2831.6223
$ws.Range("H61").Value = 1042.5385
$ws.Range("I61").Value = 850.1053000000001
$ws.Range("J61").Value = 1564.8572
$ws.Range("K61").Value = 850.1053000000001
$ws.Range("L61").Value = 1564.8572
$ws.Range("M61").Value = -638.1053000000001
$ws.Range("N61").Value = -1988.8572
$ws.Range("H74").Value = 903.913
$ws.Range("I74").Value = 881.4737
$ws.Range("K74").Value = 881.4737
$ws.Range("M74").Value = -7.473700000000008
$ws.Range("H77").Value = 903.913
$ws.Range("I77").Value = 881.4737
$ws.Range("K77").Value = 4407.3685
$ws.Range("M77").Value = -39.36850000000049
$ws.Range("H132").Value = 2769.4736
$ws.Range("I132").Value = 2504.8965
$ws.Range("J132").Value = 3622
$ws.Range("K132").Value = 7514.689499999999
$ws.Range("L132").Value = 10866
$ws.Range("M132").Value = -4984.689499999999
$ws.Range("N132").Value = -15926
$ws.Range("H136").Value = 1042.5385
$ws.Range("I136").Value = 850.1053000000001
$ws.Range("J136").Value = 1564.8572
$ws.Range("K136").Value = 2550.3159
$ws.Range("L136").Value = 4694.571599999999
$ws.Range("M136").Value = -0.315900000000056
$ws.Range("N136").Value = -9794.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
$ws.Range("H99").Value = 29413184
$ws.Range("I99").Value = 35715616
$ws.Range("J99").Value = 1838.3334
$ws.Range("K99").Value = 35715616
$ws.Range("L99").Value = 1838.3334
$ws.Range("M99").Value = -35714118
$ws.Range("N99").Value = -4834.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1182.6274
$ws.Range("I31").Value = 861.0454999999999
$ws.Range("J31").Value = 1426.5862
$ws.Range("K31").Value = 861.0454999999999
$ws.Range("L31").Value = 1426.5862
$ws.Range("M31").Value = -566.0454999999999
$ws.Range("N31").Value = -2016.5862
$ws.Range("H34").Value = 1182.6274
$ws.Range("I34").Value = 861.0454999999999
$ws.Range("J34").Value = 1426.5862
$ws.Range("K34").Value = 861.0454999999999
$ws.Range("L34").Value = 1426.5862
$ws.Range("M34").Value = -659.0454999999999
$ws.Range("N34").Value = -1830.5862
$ws.Range("H58").Value = 958.175
$ws.Range("I58").Value = 942.8077
$ws.Range("J58").Value = 986.7143
$ws.Range("K58").Value = 942.8077
$ws.Range("L58").Value = 986.7143
$ws.Range("M58").Value = -739.8077
$ws.Range("N58").Value = -1392.7143
$ws.Range("H62").Value = 4350345.5
$ws.Range("I62").Value = 2575.6445
$ws.Range("K62").Value = 2575.6445
$ws.Range("M62").Value = -1951.6445
$ws.Range("H65").Value = 4350345.5
$ws.Range("I65").Value = 2575.6445
$ws.Range("K65").Value = 12878.2225
$ws.Range("M65").Value = -9758.2225
$ws.Range("H99").Value = 1622.6364
$ws.Range("I99").Value = 1847.2
$ws.Range("J99").Value = 1435.5
$ws.Range("K99").Value = 1847.2
$ws.Range("L99").Value = 1435.5
$ws.Range("M99").Value = -349.2
$ws.Range("N99").Value = -4431.5
$ws.Range("H107").Value = 554.65216
$ws.Range("I107").Value = 238.85715
$ws.Range("J107").Value = 692.8125
$ws.Range("K107").Value = 238.85715
$ws.Range("L107").Value = 692.8125
$ws.Range("M107").Value = 1681.14285
$ws.Range("N107").Value = -4532.8125
$ws.Range("H111").Value = 3702
$ws.Range("J111").Value = 3702
$ws.Range("L111").Value = 3702
$ws.Range("N111").Value = -11882
$ws.Range("H126").Value = 1622.6364
$ws.Range("I126").Value = 1847.2
$ws.Range("J126").Value = 1435.5
$ws.Range("K126").Value = 5541.6
$ws.Range("L126").Value = 4306.5
$ws.Range("M126").Value = -3071.6
$ws.Range("N126").Value = -9246.5
$ws.Range("H132").Value = 3427.4443
$ws.Range("I132").Value = 3633.238
$ws.Range("J132").Value = 2707.1667
$ws.Range("K132").Value = 10899.714
$ws.Range("L132").Value = 8121.500100000001
$ws.Range("M132").Value = -8369.714
$ws.Range("N132").Value = -13181.5001
$ws.Range("H134").Value = 1360.6571
$ws.Range("I134").Value = 1338.4166
$ws.Range("K134").Value = 4015.2498
$ws.Range("M134").Value = -1480.2498
$ws.Range("H136").Value = 958.175
$ws.Range("I136").Value = 942.8077
$ws.Range("J136").Value = 986.7143
$ws.Range("K136").Value = 2828.4231
$ws.Range("L136").Value = 2960.1429
$ws.Range("M136").Value = -278.4231
$ws.Range("N136").Value = -8060.1429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1673.0217
$ws.Range("I68").Value = 748.25
$ws.Range("J68").Value = 2166.2334
$ws.Range("K68").Value = 2244.75
$ws.Range("L68").Value = 6498.7002
$ws.Range("M68").Value = -1433.75
$ws.Range("N68").Value = -8120.7002
$ws.Range("H71").Value = 1673.0217
$ws.Range("I71").Value = 748.25
$ws.Range("J71").Value = 2166.2334
$ws.Range("K71").Value = 6734.25
$ws.Range("L71").Value = 19496.1006
$ws.Range("M71").Value = -2678.25
$ws.Range("N71").Value = -27608.1006
$ws.Range("H131").Value = 27028602
$ws.Range("J131").Value = 1749.7667
$ws.Range("L131").Value = 5249.300099999999
$ws.Range("N131").Value = -15329.3001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6017167
$ws.Range("J11").Value = 3004
$ws.Range("L11").Value = 3004
$ws.Range("N11").Value = -3282

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1358.4667
$ws.Range("I16").Value = 992.5
$ws.Range("J16").Value = 2090.4
$ws.Range("K16").Value = 992.5
$ws.Range("L16").Value = 2090.4
$ws.Range("M16").Value = -822.5
$ws.Range("N16").Value = -2430.4
$ws.Range("H22").Value = 1074.9166
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1189.9
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 1189.9
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -1779.9
$ws.Range("H27").Value = 1074.9166
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 1189.9
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 1189.9
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -1403.9
$ws.Range("H75").Value = 12000
$ws.Range("J75").Value = 12000
$ws.Range("L75").Value = 12000
$ws.Range("N75").Value = -13872
$ws.Range("H78").Value = 12000
$ws.Range("J78").Value = 12000
$ws.Range("L78").Value = 36000
$ws.Range("N78").Value = -45360
$ws.Range("H132").Value = 28752.027
$ws.Range("I132").Value = 1392.8077
$ws.Range("J132").Value = 93419.27
$ws.Range("K132").Value = 4178.4231
$ws.Range("L132").Value = 280257.81
$ws.Range("M132").Value = -1648.4231
$ws.Range("N132").Value = -285317.81

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 428.12
$ws.Range("I136").Value = 352.75
$ws.Range("K136").Value = 1058.25
$ws.Range("M136").Value = 1491.75
